$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '40.218.10'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.21%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.212.75'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.55%  '

$ws.Range("E4").Value = '  -0.02%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '296.52'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.67%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '87.74'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("E7").Value = '  +0.48%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  -0.36%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '52.49'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +7.69%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '30.90'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.49%  '

$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("E13").Value = '  +2.25%  '

$ws.Range("E14").Value = '  -1.04%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.556.48'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.59%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '13.86'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.41%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.212.56'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.30%  '

$ws.Range("E18").Value = '  +0.92%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '40.119.10'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.11%  '

$ws.Range("E20").Value = '  -0.07%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '11.33'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.30%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.77'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.88%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '65.74'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.10%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '235.72'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.43%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("E26").Value = '  +1.32%  '

$ws.Range("E27").Value = '  -1.11%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '23.27'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.32%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.07'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -5.20%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '156.33'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.10%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '32.15'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.97%  '

$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("E34").Value = '  +0.25%  '

$ws.Range("E35").Value = '  +3.58%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0715'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.65%  '

$ws.Range("E37").Value = '  -0.64%  '

$ws.Range("E38").Value = '  +1.69%  '

$ws.Range("E39").Value = '  +2.66%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.73'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.24%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '15.57'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.10%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '3.82'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.25%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.066.67'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.07%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '19.19'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +4.21%  '

$ws.Range("E45").Value = '  +0.78%  '

$ws.Range("E46").Value = '  +0.79%  '

$ws.Range("E47").Value = '  +5.88%  '

$ws.Range("E48").Value = '  -10.98%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.429.36'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.35%  '

$ws.Range("E50").Value = '  +2.21%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.46'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.35%  '
